# REQ-51..53 System test cases added to "Test Cases & Results" (rows 54-56).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases & Results")

# --- Row 54 (TC 52 / REQ-51): Flask web app availability for library staff ---
$ws.Range("E54").Value = 'Mid Impact'
$ws.Range("F54").Value = 'Test that A Website running on flask should be available for library staff to access'
$ws.Range("G54").Value = 'WebApp/app.py is ran seperately from main.py in Src'
$ws.Range("H54").Value = 'On a browser while connecting to the same network as the RPI, type in the RPI''s IP address:5000 '
$ws.Range("I54").Value = 'A mainpage website shows'
$ws.Range("J54").Value = 'A mainpage website shows'

# --- Row 55 (TC 53 / REQ-52): Book search function on the website ---
$ws.Range("E55").Value = 'Mid Impact'
$ws.Range("F55").Value = 'Test that The website should have a search function for library staff to search a book by its title and see its details like onloan status, reserved status, borrower id, date of borrowing, bookid and location'
$ws.Range("G55").Value = 'WebApp/app.py is ran seperately from main.py in Src'
$ws.Range("H55").Value = 'On the mainmenu press the go to book search button, if a text is typed in the text field and the search button is pressed, the list of books will be narrowed to the books which contain the text in the title'
$ws.Range("I55").Value = 'A list of filtered books shows'
$ws.Range("J55").Value = 'A list of filtered books shows'

# --- Row 56 (TC 54 / REQ-53): Filter book list by location/reservation/loan status ---
$ws.Range("E56").Value = 'Mid Impact'
$ws.Range("F56").Value = 'Test that The Staff should also have the option to filter the list of books by location, reservation status and loan status'
$ws.Range("G56").Value = 'WebApp/app.py is ran seperately from main.py in Src'
$ws.Range("H56").Value = 'Same steps as REQ-52 but now with text fields for location, reservation status and loan status'
$ws.Range("I56").Value = 'A list of filtered books shows'
$ws.Range("J56").Value = 'A list of filtered books shows'

# Rows wrap their new text, so Excel grew them to fit (matches authored heights).
$ws.Rows.Item(54).RowHeight = 57.6
$ws.Rows.Item(55).RowHeight = 115.2
$ws.Rows.Item(56).RowHeight = 57.6

# Move the live selection to where the author left off editing.
$null = $ws.Range("I58").Select()
